$wb = $excel.ActiveWorkbook
$wsVisuais = $wb.Worksheets.Item("Visuais")
$wsMedidas = $wb.Worksheets.Item("Medidas")

# --- Sheet "Visuais": add newly documented visual rows 24-26 (columns D:F) ---
$wsVisuais.Range("D24").Value = "Ano"
$wsVisuais.Range("E24").Value = "Campo"
$wsVisuais.Range("F24").Value = "fat-partidas"

$wsVisuais.Range("D25").Value = "Partidas"
$wsVisuais.Range("E25").Value = "Medida"
$wsVisuais.Range("F25").Value = "Medidas"

$wsVisuais.Range("D26").Value = "Seleções"
$wsVisuais.Range("E26").Value = "Medida"
$wsVisuais.Range("F26").Value = "Medidas"

# --- Sheet "Medidas": formatting updates ---
# Turn on wrap text for the "Medida"/measure rows in column A (A5:A16)
$wsMedidas.Range("A5:A16").WrapText = $true

# Narrow columns B, C, D so the long description/calculation text wraps onto more lines
$wsMedidas.Columns.Item(2).ColumnWidth = 10.5
$wsMedidas.Columns.Item(3).ColumnWidth = 22.5
$wsMedidas.Columns.Item(4).ColumnWidth = 37

# Explicit row heights reflecting the new wrapped-text layout
$wsMedidas.Rows.Item(2).RowHeight = 30
$wsMedidas.Rows.Item(3).RowHeight = 45
$wsMedidas.Rows.Item(6).RowHeight = 30
$wsMedidas.Rows.Item(8).RowHeight = 30
$wsMedidas.Rows.Item(9).RowHeight = 45
$wsMedidas.Rows.Item(10).RowHeight = 30
$wsMedidas.Rows.Item(11).RowHeight = 105
$wsMedidas.Rows.Item(12).RowHeight = 30
$wsMedidas.Rows.Item(13).RowHeight = 30
$wsMedidas.Rows.Item(14).RowHeight = 60
$wsMedidas.Rows.Item(15).RowHeight = 45
$wsMedidas.Rows.Item(17).RowHeight = 180

# Update the selection/scroll position left behind on "Medidas"
$wsMedidas.Activate()
$wsMedidas.Range("A1:D17").Select()

# Switch the active tab back to "Visuais", like in the saved workbook
$wsVisuais.Activate()
$wsVisuais.Range("A1").Select()
